# Update the weekly "Hortaliza, Terminal La Palmera de La Serena - Alcachofa"
# price report: dates and associated price/origin data were refreshed
# (fruta / hortaliza, semanal).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 44484
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 9500
$ws.Range("P2").Value = 317

# Row 3
$ws.Range("D3").Value = 44420
$ws.Range("H3").Value = 'Madrigal'
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14500
$ws.Range("N3").Value = '$/caja 40 unidades'
$ws.Range("O3").Value = 'Provincia de Limarí'
$ws.Range("P3").Value = 362
$ws.Range("Q3").Value = 40

# Row 4
$ws.Range("D4").Value = 44420
$ws.Range("H4").Value = 'Madrigal'
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("N4").Value = '$/caja 40 unidades'
$ws.Range("O4").Value = 'Provincia del Elquí'
$ws.Range("P4").Value = 338
$ws.Range("Q4").Value = 40

# Row 5
$ws.Range("D5").Value = 44438
$ws.Range("H5").Value = 'Española'
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 11500
$ws.Range("N5").Value = '$/caja 30 unidades'
$ws.Range("O5").Value = 'Provincia del Elquí'
$ws.Range("P5").Value = 383
$ws.Range("Q5").Value = 30

# Row 6
$ws.Range("D6").Value = 44498
$ws.Range("H6").Value = 'Española'
$ws.Range("K6").Value = 8500
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 8750
$ws.Range("N6").Value = '$/caja 30 unidades'
$ws.Range("P6").Value = 292
$ws.Range("Q6").Value = 30

# Row 7
$ws.Range("D7").Value = 44426
$ws.Range("J7").Value = 600
$ws.Range("K7").Value = 11500
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 11750
$ws.Range("P7").Value = 392

# Row 8
$ws.Range("D8").Value = 44426
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 12500
$ws.Range("L8").Value = 13000
$ws.Range("M8").Value = 12750
$ws.Range("P8").Value = 319

# Row 9
$ws.Range("D9").Value = 44427
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 12500
$ws.Range("O9").Value = 'Provincia de Limarí'
$ws.Range("P9").Value = 312
